$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# A new "after_HSA" column is inserted between the existing HSA (I) and
# total_savings (J) columns. The original "total_savings" header/data
# slide over to the new column K, and J becomes "after_HSA" with a new
# set of (smaller) values.
$ws.Range("J1").Value = "after_HSA"
$ws.Range("K1").Value = "total_savings"

# --- Data rows ----------------------------------------------------------
# Columns: A=coveragetier B=xt_totpredrisk_l1 C=oop_IYC_model D=premium_HD
#          E=IYC_Total_Cost F=oop_HDHP_model G=premium_HDHP H=HD_Total_Cost
#          I=HSA J=after_HSA(new values) K=total_savings(old J values, rounded)
$data = @(
    @{ Row = 2;  C = 410;  E = 1994; F = 670;  H = 1258; J = 406;  K = 1588 },
    @{ Row = 3;  C = 460;  E = 2044; F = 810;  H = 1398; J = 546;  K = 1498 },
    @{ Row = 4;  C = 510;  E = 2094; F = 970;  H = 1558; J = 706;  K = 1388 },
    @{ Row = 5;  C = 600;  E = 2184; F = 1250; H = 1838; J = 986;  K = 1198 },
    @{ Row = 6;  C = 970;  E = 2554; F = 2180; H = 2768; J = 1916; K = 638  },
    @{ Row = 7;  C = 1040; E = 4988; F = 2470; H = 3934; J = 2230; K = 2758 },
    @{ Row = 8;  C = 1170; E = 5118; F = 2800; H = 4264; J = 2560; K = 2558 },
    @{ Row = 9;  C = 1290; E = 5238; F = 3080; H = 4544; J = 2840; K = 2398 },
    @{ Row = 10; C = 1460; E = 5408; F = 3490; H = 4954; J = 3250; K = 2158 },
    @{ Row = 11; C = 2090; E = 6038; F = 4720; H = 6184; J = 4480; K = 1558 }
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Range("C$r").Value = $d.C
    $ws.Range("E$r").Value = $d.E
    $ws.Range("F$r").Value = $d.F
    $ws.Range("H$r").Value = $d.H
    $ws.Range("J$r").Value = $d.J
    $ws.Range("K$r").Value = $d.K
}

# Apply the same integer number format ("0") already used by the other
# numeric columns (B-J, style index 1) to the newly-added column K.
$ws.Range("K2:K11").NumberFormat = "0"
